# Automatische test-sync: 2025-06-22 21:37:50
# Adds a new "Open sollicitatie" log entry to the Logs sheet (row 42),
# extends the conditional-formatting ranges to cover it, and bumps the
# matching Dashboard summary count for "Sollicitatie / Vacature".

$wb = $excel.ActiveWorkbook

$logs = $wb.Worksheets.Item("Logs")

$newRow = 42
$logs.Cells.Item($newRow, 1).Value = "Open sollicitatie"
$logs.Cells.Item($newRow, 2).Value = "mailmind.test@zohomail.eu"
$logs.Cells.Item($newRow, 3).Value = "Zijn er op dit moment openstaande functies bij jullie bedrijf?"
$logs.Cells.Item($newRow, 4).Value = "Sollicitatie / Vacature"
$logs.Cells.Item($newRow, 6).Value = "2025-06-22 21:37:13"
$logs.Cells.Item($newRow, 7).Value = "Nee"

# Extend the conditional formatting ranges (D2:D41 -> D2:D42, G2:G41 -> G2:G42)
# so the new row is covered by the existing colour rules.
$dRule = $logs.Range("D2:D41").FormatConditions.Item(1)
$dRule.ModifyAppliesToRange($logs.Range("D2:D42"))

$gRule = $logs.Range("G2:G41").FormatConditions.Item(1)
$gRule.ModifyAppliesToRange($logs.Range("G2:G42"))

# Dashboard: "Sollicitatie / Vacature" count 4 -> 5
$dashboard = $wb.Worksheets.Item("Dashboard")
$dashboard.Cells.Item(3, 2).Value = 5
